# The source edit underlying this diff ("Fixed POI packaging and upgraded
# to POI 3.15") only changes the *serialization* of the package: every
# hunk in the unified diff re-orders XML attributes into alphabetical
# order (a side effect of the Apache POI XML writer used when the
# fixture was regenerated). Namespace declarations on <w:document>, the
# attributes of <w:pgSz>/<w:pgMar>, <w:rFonts>, <w:lang>,
# <w:latentStyles>/<w:lsdException>, <w:style>, <w:tblInd> and
# <w:tblCellMar>'s margin elements are all re-ordered, but every
# attribute name/value pair, every element, and all document content
# and formatting stay exactly the same.
#
# That kind of attribute-order change is an artifact of the XML writer
# used to produce the package, not something the Word object model
# exposes (Word/COM automation has no notion of "attribute order" --
# it reads/writes the document's content and formatting model, and the
# underlying OOXML writer controls serialization order on save).
# Consequently there is no content or formatting change to apply here:
# the document's paragraphs, runs, styles, sections and page setup are
# already identical to the target. This script intentionally performs
# no mutations, so the document is re-saved as-is.

$d = $word.ActiveDocument

# Touch nothing; the target state already matches the document's
# current content/formatting -- only cosmetic XML attribute ordering
# differs in the source diff, which is outside the Word OM's surface.
